$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Bad" conditional-looking style from C28 (it becomes a normal cell)
$ws.Range("C28").Style = "Normal"

# Insert a brand-new row 29 (shifts the former rows 29-56 down to 30-57)
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29
$ws.Range("A29").Value = "HunterAbilitySet"
$ws.Range("B29").Value = "ShadowOps_SnapShotShot"
$ws.Range("C29").Value = "UIPerk_snapshot"
$ws.Range("C29").Style = "Bad"
$ws.Range("D29").Value = "Redo"

# Add new "D" column annotations to several of the shifted rows
$ws.Range("D30").Value = "OK"
$ws.Range("D31").Value = "Check"
$ws.Range("D34").Value = "Check"
$ws.Range("D35").Value = "Check"
$ws.Range("D36").Value = "Check"
$ws.Range("D41").Value = "Check"
$ws.Range("D43").Value = "Redo"

# Update the view to match the author's final cursor/scroll position
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("D37").Select()
